$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 10 data first (so new shared strings are interned in the
# same order the original author's edit produced them)
$ws.Range("A10").Value = "Quick Fit Terminal, Male .250"""
$ws.Range("B10").Value = 1
$ws.Range("C10").Value = "1287-R"
$ws.Range("D10").Value = "N/A"
$ws.Range("E10").Value = "N/A"

# Fill in MPN (column C) for existing rows that previously had no MPN
$ws.Range("C4").Value = "1N4148W"
$ws.Range("C5").Value = "352210MJT"

# Update sheet view: zoom to 100% (normal), and change selected cell
$ws.Application.ActiveWindow.Zoom = 100
$ws.Range("D19").Select()
